# NSMB - 8-8 & some of 8-F2 done.
# Append rows 192-200 to the "V4" sheet, extending the existing checkpoint
# / split-time log table (columns A..D) that previously ended at row 191.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# row, Name (col A), Start tick (col B), End tick (col C).
# Col D keeps the existing table formula: =IF(Bn > 0,Cn-Bn, 0)
$rows = @(
    @(192, "Checkpoint 1633",       59432, 69850),
    @(193, "Checkpoint 1946/1944",  59610, 70031),
    @(194, "Checkpoint 2388/2385",  59716, 70137),
    @(195, "Checkpoint 2896/2895",  59863, 70285),
    @(196, "Get flag",              60022, 70444),
    @(197, "End Level",             60540, 70962),
    @(198, "Enter 8-F2",            60909, 71707),
    @(199, "1st Move",              61137, 71957),
    @(200, "Platform 1st Move",     61238, 72069)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Formula = "=IF(B$rowNum > 0,C$rowNum-B$rowNum, 0)"
}

# Move the active selection to the next blank row below the appended data,
# matching where Excel leaves the cursor after typing in the last new row.
$ws.Activate()
$ws.Range("B201").Select()
